# regen sval data to filter save games
# Update the B:G numeric columns (TB, d2S, K, IP, Win, sum) for rows 2-17
# on the active worksheet with the newly regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.0001488876196638067, 0.002777888934908601, 9844.520545567508, 616238.5361209477, 0, 626083.0595932917)
    3  = @(0.6753301551942219, 1.667794583268128, 9844.520545567508, 645.3272768299601, 1, 10492.19094713593)
    4  = @(1.459612070389937, 10.29869402782916, 3.900430680208489, 616238.5361209477, 0, 616254.1948577261)
    5  = @(0.003994804209775715, 0.002777888934908601, 26.21740644021617, 8.660232485948974, 0, 34.88441161930983)
    6  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    7  = @(0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 645.3272768299601, 0, 646.1903192409343)
    8  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    9  = @(0.04763786555579896, 1.667794583268128, 3.900430680208489, 616238.5361209477, 0, 616244.1519840767)
    10 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    11 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    12 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 1, 13.71653804550039)
    13 = @(0.127881588408715, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0, 0.8245905782990299)
    14 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 0, 17.45944343273191)
    15 = @(0.003994804209775715, 3099.503889238888, 0.8054896365839992, 645.3272768299601, 1, 3745.640650509642)
    16 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 645.3272768299601, 1, 674.6720899238343)
    17 = @(3.230985683306322, 3099.503889238888, 0.8054896365839992, 8.660232485948974, 1, 3112.200597044728)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
